# DP Compress the code
#
# The "software location" path stored in B2 is rewritten from the old,
# developer-local Windows path to the shorter network-share / git-repo
# path (and forward slashes), e.g.:
#   C:\Users\dpere\Documents\JTMT\forecast\create_forecast_basic\current
#   -> W:/Data/Forecast/Tools/forecast_git/create_forecast_basic/current
#
# The author also re-confirmed the cell format (Format Cells > OK) for the
# whole A1:B3 block / columns A:B while editing, which is why the saved
# workbook shows every cell/column pointing at an explicitly-applied
# (but visually unchanged: General horizontal, Bottom vertical, locked,
# not hidden) style instead of the original implicit default one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- the actual content edit -------------------------------------------
$ws.Range("B2").Value = "W:/Data/Forecast/Tools/forecast_git/create_forecast_basic/current"

# --- re-apply (unchanged) formatting to the used range + its columns ---
$xlHAlignGeneral = 1
$xlVAlignBottom = -4107

foreach ($target in @($ws.Range("A1:B3"), $ws.Range("A1:B1").EntireColumn)) {
    $target.HorizontalAlignment = $xlHAlignGeneral
    $target.VerticalAlignment = $xlVAlignBottom
    $target.WrapText = $false
    $target.IndentLevel = 0
    $target.Locked = $true
    $target.FormulaHidden = $false
}
